$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updates to column B (predicted/ASR token) and column C (count) per row.
# Row => [B value, C value]
$updates = @{
    2  = @("<that>", 30)
    3  = @($null, 27)
    4  = @("<not>", 37)
    5  = @($null, 29)
    6  = @($null, 32)
    7  = @($null, 33)
    8  = @("<nun>", 32)
    9  = @($null, 32)
    10 = @("<past>", $null)
    11 = @("<ban>", 32)
    12 = @($null, 32)
    13 = @($null, 35)
    14 = @($null, 34)
    16 = @("<number>", $null)
    17 = @("<encape>", 35)
    18 = @("<left>", 25)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $bVal = $vals[0]
    $cVal = $vals[1]

    if ($null -ne $bVal) {
        $ws.Cells.Item($row, 2).Value = $bVal
    }
    if ($null -ne $cVal) {
        $ws.Cells.Item($row, 3).Value = $cVal
    }
}
